$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, value first
$ws.Range("H1").Value = "Save"

# Copy the formatting from the existing header cell (G1) onto H1 so the new
# header matches the look of the other header cells (bold, centered, bordered)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add new value in H2 for the corresponding data row
$ws.Range("H2").Value = 1
